$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 667.8333
$ws.Range("I16").Value = 254
$ws.Range("J16").Value = 874.75
$ws.Range("K16").Value = 254
$ws.Range("L16").Value = 874.75
$ws.Range("M16").Value = -24
$ws.Range("N16").Value = -1334.75

$ws.Range("H40").Value = 2000
$ws.Range("J40").Value = 2000
$ws.Range("L40").Value = 2000
$ws.Range("N40").Value = -2350

$ws.Range("H43").Value = 13808
$ws.Range("J43").Value = 13808
$ws.Range("L43").Value = 13808
$ws.Range("N43").Value = -13946

$ws.Range("H49").Value = 1406.9166
$ws.Range("I49").Value = 270.42856
$ws.Range("J49").Value = 2998
$ws.Range("K49").Value = 811.28568
$ws.Range("L49").Value = 8994
$ws.Range("M49").Value = -675.28568
$ws.Range("N49").Value = -9266

$ws.Range("H98").Value = 3180.524
$ws.Range("J98").Value = 3496
$ws.Range("L98").Value = 3496
$ws.Range("N98").Value = -6492

$ws.Range("H122").Value = 3180.524
$ws.Range("J122").Value = 3496
$ws.Range("L122").Value = 10488
$ws.Range("N122").Value = -15388

$ws.Range("H136").Value = 84618.46000000001
$ws.Range("J136").Value = 84618.46000000001
$ws.Range("L136").Value = 84618.46000000001
$ws.Range("N136").Value = -94818.46000000001

$ws.Range("H138").Value = 2399.04
$ws.Range("I138").Value = 3316.5386
$ws.Range("J138").Value = 2076.6758
$ws.Range("K138").Value = 9949.6158
$ws.Range("L138").Value = 6230.0274
$ws.Range("M138").Value = -4809.6158
$ws.Range("N138").Value = -16510.0274

$ws.Range("H141").Value = 3907.6875
$ws.Range("I141").Value = 3586.182
$ws.Range("K141").Value = 10758.546
$ws.Range("M141").Value = -5578.545999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1614376.9
$ws.Range("J5").Value = 93.55556
$ws.Range("L5").Value = 93.55556
$ws.Range("N5").Value = -317.55556

$ws.Range("H32").Value = 4388638.5
$ws.Range("I32").Value = 4903517.5
$ws.Range("J32").Value = 12166.167
$ws.Range("K32").Value = 4903517.5
$ws.Range("L32").Value = 12166.167
$ws.Range("M32").Value = -4903230.5
$ws.Range("N32").Value = -12740.167

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""

$ws.Range("H110").Value = 8269952
$ws.Range("I110").Value = 10107575
$ws.Range("K110").Value = 10107575
$ws.Range("M110").Value = -10105530

$ws.Range("H132").Value = 2135333.5
$ws.Range("I132").Value = 1205.4524
$ws.Range("K132").Value = 3616.357199999999
$ws.Range("M132").Value = -1086.357199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1614376.9
$ws.Range("J4").Value = 93.55556
$ws.Range("L4").Value = 93.55556
$ws.Range("N4").Value = -323.55556

$ws.Range("H6").Value = 69171.5
$ws.Range("I6").Value = 38647
$ws.Range("K6").Value = 38647
$ws.Range("M6").Value = -38534

$ws.Range("H22").Value = 23811690
$ws.Range("I22").Value = 23811690
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 23811690
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -23811517
$ws.Range("N22").Value = ""

$ws.Range("H68").Value = 70000
$ws.Range("J68").Value = 70000
$ws.Range("L68").Value = 70000
$ws.Range("N68").Value = -71622

$ws.Range("H71").Value = 70000
$ws.Range("J71").Value = 70000
$ws.Range("L71").Value = 210000
$ws.Range("N71").Value = -218112

$ws.Range("H94").Value = 1680.4468
$ws.Range("I94").Value = 1352.6666
$ws.Range("J94").Value = 2453.0715
$ws.Range("K94").Value = 1352.6666
$ws.Range("L94").Value = 2453.0715
$ws.Range("M94").Value = -901.6666
$ws.Range("N94").Value = -3355.0715

$ws.Range("H99").Value = 10035.233
$ws.Range("I99").Value = 10867.615
$ws.Range("K99").Value = 10867.615
$ws.Range("M99").Value = -9369.615

$ws.Range("H134").Value = 40565.727
$ws.Range("I134").Value = 43592.207
$ws.Range("K134").Value = 130776.621
$ws.Range("M134").Value = -128241.621

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20685.812
$ws.Range("I31").Value = 999.6667
$ws.Range("K31").Value = 999.6667
$ws.Range("M31").Value = -704.6667

$ws.Range("H34").Value = 20685.812
$ws.Range("I34").Value = 999.6667
$ws.Range("K34").Value = 999.6667
$ws.Range("M34").Value = -797.6667

$ws.Range("H132").Value = 23811956
$ws.Range("I132").Value = 2492.2778
$ws.Range("K132").Value = 7476.8334
$ws.Range("M132").Value = -4946.8334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 146.61539
$ws.Range("I2").Value = 183.33333
$ws.Range("J2").Value = 64
$ws.Range("K2").Value = 1099.99998
$ws.Range("L2").Value = 384
$ws.Range("M2").Value = -986.9999800000001
$ws.Range("N2").Value = -610

$ws.Range("H9").Value = 129028.57
$ws.Range("I9").Value = 200149.5
$ws.Range("J9").Value = 100580.2
$ws.Range("K9").Value = 600448.5
$ws.Range("L9").Value = 301740.6
$ws.Range("M9").Value = -600224.5
$ws.Range("N9").Value = -302188.6

$ws.Range("H22").Value = 5327.5713
$ws.Range("J22").Value = 4709.4
$ws.Range("L22").Value = 14128.2
$ws.Range("N22").Value = -14466.2

$ws.Range("H27").Value = 5327.5713
$ws.Range("J27").Value = 4709.4
$ws.Range("L27").Value = 14128.2
$ws.Range("N27").Value = -14332.2

$ws.Range("H131").Value = 1437.8788
$ws.Range("J131").Value = 1480.1183
$ws.Range("L131").Value = 4440.3549
$ws.Range("N131").Value = -14520.3549

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 2375
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 8000
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 8000
$ws.Range("M22").Value = 29
$ws.Range("N22").Value = -9058

$ws.Range("H102").Value = 12290447
$ws.Range("I102").Value = 12290447
$ws.Range("K102").Value = 12290447
$ws.Range("M102").Value = -12288825

$ws.Range("H122").Value = 1700272.5
$ws.Range("I122").Value = 1999496.8
$ws.Range("K122").Value = 5998490.4
$ws.Range("M122").Value = -5996040.4

$ws.Range("H126").Value = 9767104
$ws.Range("I126").Value = 3600957.8
$ws.Range("K126").Value = 10802873.4
$ws.Range("M126").Value = -10800403.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 172799.8
$ws.Range("I25").Value = 13999
$ws.Range("J25").Value = 212500
$ws.Range("K25").Value = 13999
$ws.Range("L25").Value = 212500
$ws.Range("M25").Value = -13769
$ws.Range("N25").Value = -212960

$ws.Range("H40").Value = 5887136
$ws.Range("I40").Value = 5100
$ws.Range("J40").Value = 8408009
$ws.Range("K40").Value = 5100
$ws.Range("L40").Value = 8408009
$ws.Range("M40").Value = -4964
$ws.Range("N40").Value = -8408281

$ws.Range("H55").Value = 1618.3462
$ws.Range("J55").Value = 1837.9231
$ws.Range("L55").Value = 1837.9231
$ws.Range("N55").Value = -2183.9231

$ws.Range("H122").Value = 41983550
$ws.Range("I122").Value = 49597416
$ws.Range("K122").Value = 148792248
$ws.Range("M122").Value = -148789798

$ws.Range("H132").Value = 939933.9399999999
$ws.Range("I132").Value = 2471.1785
$ws.Range("K132").Value = 7413.5355
$ws.Range("M132").Value = -4883.5355

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 500
$ws.Range("K5").Value = 500
$ws.Range("M5").Value = -388

$ws.Range("H15").Value = 60783730
$ws.Range("J15").Value = 60783730
$ws.Range("L15").Value = 60783730
$ws.Range("N15").Value = -60784306

$ws.Range("H86").Value = 5057499.5
$ws.Range("J86").Value = 5057499.5
$ws.Range("L86").Value = 5057499.5
$ws.Range("N86").Value = -5059745.5

$ws.Range("H89").Value = 5057499.5
$ws.Range("J89").Value = 5057499.5
$ws.Range("L89").Value = 25287497.5
$ws.Range("N89").Value = -25298729.5

$ws.Range("H107").Value = 1140.1052
$ws.Range("I107").Value = 1232.0625
$ws.Range("K107").Value = 3696.1875
$ws.Range("M107").Value = -1776.1875

$ws.Range("H113").Value = 3188.261
$ws.Range("I113").Value = 3037.7273
$ws.Range("J113").Value = 6500
$ws.Range("K113").Value = 9113.1819
$ws.Range("L113").Value = 19500
$ws.Range("M113").Value = -6943.1819
$ws.Range("N113").Value = -23840

$ws.Range("H132").Value = 4298.811
$ws.Range("I132").Value = 1677.4849
$ws.Range("K132").Value = 5032.4547
$ws.Range("M132").Value = -2502.4547
